$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (the FAPs -> Resolving-Mac record) entirely.
$ws.Rows("5:5").Delete()

# Update row 2 values (Target cluster = ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05003266666666667
$ws.Range("H2").Value = 0.150098
$ws.Range("M2").Value = 30.58864766666666
$ws.Range("N2").Value = 91.76594299999999
$ws.Range("O2").Value = 0.3925391465174898
$ws.Range("P2").Value = 0.3925391465174898
$ws.Range("Q2").Value = 1.530431612490444
$ws.Range("R2").Value = 13.773884512414
$ws.Range("S2").Value = 0.3925391465174898
$ws.Range("T2").Value = 0.3925391465174898

# Update row 3 values (Target cluster = FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05003266666666667
$ws.Range("H3").Value = 0.150098
$ws.Range("O3").Value = 0.291183949679193
$ws.Range("P3").Value = 0.291183949679193
$ws.Range("Q3").Value = 1.135267974143333
$ws.Range("R3").Value = 10.21741176729
$ws.Range("S3").Value = 0.291183949679193
$ws.Range("T3").Value = 0.291183949679193

# Update row 4 values (Target cluster = MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05003266666666667
$ws.Range("H4").Value = 0.150098
$ws.Range("M4").Value = 24.64590566666666
$ws.Range("N4").Value = 73.93771699999999
$ws.Range("O4").Value = 0.3162769038033173
$ws.Range("P4").Value = 0.3162769038033172
$ws.Range("Q4").Value = 1.233100382918444
$ws.Range("R4").Value = 11.097903446266
$ws.Range("S4").Value = 0.3162769038033173
$ws.Range("T4").Value = 0.3162769038033172
